# Update share price data on Sheet1 (changed share names/values in cash).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row => B,C,D,E,F,G values (HIGH, LOW, CLOSE, LTP, VOL, 9:25 CLOSE)
$data = @{
    2  = @(1480.65, 1469.25, 3028.9,   1470.9,   20,  1477.75)
    3  = @(1567.5,  1538,    1540.25,  1542.1,   15,  0)
    4  = @(51775,   51641.05,51712,    51727.6,  8,   51736.3)
    5  = @(109.3,   108.09,  108.77,   108.8,    423, 108.95)
    6  = @(854.5,   842.65,  844,      844.1,    25,  853)
    7  = @(674.7,   670.65,  671.6,    672.2,    81,  672.25)
    8  = @(1245.4,  1237.75, 1241.55,  1241.6,   161, 1244.8)
    9  = @(960.9,   947.8,   952.05,   953.9,    12,  955)
    10 = @(25313.45,25225,   25235,    25236.75, 39,  25305)
    11 = @(3049.35, 2989.1,  3004,     2999.45,  146, 3046.6)
    12 = @(824.45,  817.65,  821.15,   821.45,   167, 818.4)
    13 = @(1212.2,  1191.3,  1192.7,   1193.15,  23,  1205)
    14 = @(1081.4,  1070.1,  1071.8,   1072,     99,  1077.55)
    15 = @(153.49,  152.19,  152.57,   152.44,   439, 153.39)
    16 = @(4513.5,  4478.7,  4502.25,  4497.05,  24,  4499.65)
    17 = @(3770.95, 3633,    3728.7,   3732.6,   48,  3641)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
